$p = $ppt.ActivePresentation

# --- Slide 2 ---
$s2 = $p.Slides.Item(2)

# Shape 9 (id 218): "Aula 4| Etapa 3: " -> split into two runs:
# "Aula 4| " (unchanged run) + "Etapa 3: " (new run, same formatting)
$titleShape = $s2.Shapes.Item(9)
$titleRange = $titleShape.TextFrame.TextRange
$etapaRange = $titleRange.Characters(9, 9)
$etapaRange.Text = "Etapa 3: "

# Shape 10 (id 219): "Utilização da Classe ListaEncadeada" -> "Implementação da Classe No()"
# Also shrink the shape width.
$classShape = $s2.Shapes.Item(10)
$classShape.Width = 7923965 / 12700
$classRange = $classShape.TextFrame.TextRange
$classRange.Delete()
$classRange.Text = "Implementação da Classe No()"

# --- Slide 4 ---
$s4 = $p.Slides.Item(4)

# Shape 12 ("Imagem 13"): move picture up (change Top only).
$img = $s4.Shapes.Item(12)
$img.Top = 7.668898
